$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

$ws.Range("A1").Value = 'Datos actualizados a 23 de Junio de 2020 a las 01:52'
$ws.Range("B4").Value = 2387543
$ws.Range("C4").Value = 30886
$ws.Range("D4").Value = 991023
$ws.Range("E4").Value = 1273917
$ws.Range("G4").Value = 356
$ws.Range("H4").Value = 122603
$ws.Range("B5").Value = 1111348
$ws.Range("C5").Value = 24358
$ws.Range("E5").Value = 480715
$ws.Range("G5").Value = 748
$ws.Range("H5").Value = 51407
$ws.Range("B14").Value = 192119
$ws.Range("C14").Value = 544
$ws.Range("E14").Value = 7850
$ws.Range("A21").Value = 'Canada'
$ws.Range("B21").Value = 101637
$ws.Range("C21").Value = 300
$ws.Range("D21").Value = 64334
$ws.Range("E21").Value = 28867
$ws.Range("G21").Value = 6
$ws.Range("H21").Value = 8436
$ws.Range("A22").Value = 'Sudafrica'
$ws.Range("B22").Value = 101590
$ws.Range("C22").Value = 4288
$ws.Range("D22").Value = 53444
$ws.Range("E22").Value = 46155
$ws.Range("G22").Value = 61
$ws.Range("H22").Value = 1991
$ws.Range("B34").Value = 44931
$ws.Range("C34").Value = 2146
$ws.Range("E34").Value = 30735
$ws.Range("G34").Value = 32
$ws.Range("H34").Value = 1043
$ws.Range("B44").Value = 29157
$ws.Range("C44").Value = 324
$ws.Range("E44").Value = 19718
$ws.Range("B54").Value = 17916
$ws.Range("C54").Value = 52
$ws.Range("D54").Value = 16133
$ws.Range("E54").Value = 830
$ws.Range("B62").Value = 12772
$ws.Range("C62").Value = 466
$ws.Range("E62").Value = 11116
$ws.Range("A65").Value = 'Camerun'
$ws.Range("B65").Value = 12041
$ws.Range("D65").Value = 7740
$ws.Range("E65").Value = 3993
$ws.Range("G65").Value = 5
$ws.Range("H65").Value = 308
$ws.Range("A66").Value = 'Argelia'
$ws.Range("B66").Value = 11920
$ws.Range("C66").Value = 149
$ws.Range("D66").Value = 8559
$ws.Range("E66").Value = 2509
$ws.Range("G66").Value = 7
$ws.Range("H66").Value = 852
$ws.Range("B70").Value = 8751
$ws.Range("C70").Value = 6
$ws.Range("E70").Value = 365
$ws.Range("A85").Value = 'Gabon'
$ws.Range("B85").Value = 4739
$ws.Range("C85").Value = 311
$ws.Range("D85").Value = 2002
$ws.Range("E85").Value = 2698
$ws.Range("G85").Value = 5
$ws.Range("H85").Value = 39
$ws.Range("A86").Value = 'Etiopia'
$ws.Range("B86").Value = 4663
$ws.Range("C86").Value = 131
$ws.Range("D86").Value = 1297
$ws.Range("E86").Value = 3291
$ws.Range("G86").Value = 1
$ws.Range("H86").Value = 75
$ws.Range("A87").Value = 'Republica de Yibuti'
$ws.Range("B87").Value = 4599
$ws.Range("C87").Value = 17
$ws.Range("D87").Value = 3952
$ws.Range("E87").Value = 599
$ws.Range("G87").Value = 3
$ws.Range("H87").Value = 48
$ws.Range("B135").Value = 882
$ws.Range("C135").Value = 6
$ws.Range("D135").Value = 815
$ws.Range("E135").Value = 42
$ws.Range("B150").Value = 595
$ws.Range("C150").Value = 24
$ws.Range("D150").Value = 116
$ws.Range("E150").Value = 469
$ws.Range("D151").Value = 380
$ws.Range("E151").Value = 176
$ws.Range("D157").Value = 328
$ws.Range("E157").Value = 21
$ws.Range("D179").Value = 95
$ws.Range("E179").Value = 2
$ws.Range("A202").Value = 'Dominica'
$ws.Range("A203").Value = 'Fiyi'
